# Commit: "#5: fund, bonds, otherbonds, antique done"
#
# Sheet 5 ("具有相當價值之財產" / "property of considerable value") had its
# header row (row 1) accidentally populated with a copy of the data row
# instead of real column labels, and was missing the trailing metadata
# columns (property_category, category, date, legislator_name,
# legislator_id, source_file, index) that every other sheet in this
# workbook carries. This fixes the header and fills in the missing
# columns/values for the one data row (黃金3431克 - gold, recorded under
# the "otherbonds" property category), and also tidies the numeric-looking
# "total" value (drops the stray decimal point).

$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item(5)

# Seed the new columns (F:L) with the existing row styles (bold/bordered
# header style for row 1, plain style for row 2) before writing into them,
# so they visually match the rest of the table.
$ws5.Range("E1").Copy()
$ws5.Range("F1:L1").PasteSpecial(-4122)
$ws5.Range("E2").Copy()
$ws5.Range("F2:L2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 1: proper header labels (mirrors the land/building/car/deposit sheets).
$ws5.Range("B1").Value = "name"
$ws5.Range("C1").Value = "quantity"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "total"
$ws5.Range("F1").Value = "property_category"
$ws5.Range("G1").Value = "category"
$ws5.Range("H1").Value = "date"
$ws5.Range("I1").Value = "legislator_name"
$ws5.Range("J1").Value = "legislator_id"
$ws5.Range("K1").Value = "source_file"
$ws5.Range("L1").Value = "index"

# Row 2: the actual data, now with the full set of columns populated.
$ws5.Range("A2").Value = 89
$ws5.Range("B2").Value = "黃金3431克"
$ws5.Range("C2").Value = 1
$ws5.Range("D2").Value = "許添財"
$ws5.Range("E2").Value = "3918202(台灣銀行買進牌價1142g)"
$ws5.Range("F2").Value = "otherbonds"
$ws5.Range("G2").Value = "normal"

# Force the date column to stay a plain text value ("2013-12-31") instead
# of Excel auto-converting the typed string into a date serial number.
$ws5.Range("H2").NumberFormat = "@"
$ws5.Range("H2").Value = "2013-12-31"

$ws5.Range("I2").Value = "許添財"
$ws5.Range("J2").Value = 639
$ws5.Range("K2").Value = "tmpbb0f1"
$ws5.Range("L2").Value = 89

# Re-apply the plain row-2 style over H2 so it matches its neighbours again
# now that its value is set (the Text number format above was only needed
# transiently to stop the date auto-detection).
$ws5.Range("B2").Copy()
$ws5.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
